$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.454345666666666
$ws.Range("H2").Value = 7.363036999999999
$ws.Range("I2").Value = 0.2857469401257222
$ws.Range("J2").Value = 0.3027613480760606
$ws.Range("M2").Value = 90.43008666666667
$ws.Range("N2").Value = 271.29026
$ws.Range("O2").Value = 0.863466363695901
$ws.Range("P2").Value = 0.8656179140344247
$ws.Range("Q2").Value = 221.9466913466244
$ws.Range("R2").Value = 1997.52022211962
$ws.Range("S2").Value = 0.2467328713275877
$ws.Range("T2").Value = 0.26207564657185
$ws.Range("G3").Value = 2.454345666666666
$ws.Range("H3").Value = 7.363036999999999
$ws.Range("I3").Value = 0.2857469401257222
$ws.Range("J3").Value = 0.3027613480760606
$ws.Range("O3").Value = 0.000602137432244878
$ws.Range("P3").Value = 0.0006036378137891445
$ws.Range("Q3").Value = 0.1547743102008889
$ws.Range("R3").Value = 1.392968791808
$ws.Range("S3").Value = 0.0001720589287991333
$ws.Range("T3").Value = 0.0001827581982524874
$ws.Range("G4").Value = 2.454345666666666
$ws.Range("H4").Value = 7.363036999999999
$ws.Range("I4").Value = 0.2857469401257222
$ws.Range("J4").Value = 0.3027613480760606
$ws.Range("M4").Value = 9.467965
$ws.Range("N4").Value = 28.403895
$ws.Range("O4").Value = 0.0904043069236993
$ws.Range("P4").Value = 0.09062957269587499
$ws.Range("Q4").Value = 23.23765886990166
$ws.Range("R4").Value = 209.138929829115
$ws.Range("S4").Value = 0.02583275407763371
$ws.Range("T4").Value = 0.02743913160496044
$ws.Range("G5").Value = 2.454345666666666
$ws.Range("H5").Value = 7.363036999999999
$ws.Range("I5").Value = 0.2857469401257222
$ws.Range("J5").Value = 0.3027613480760606
$ws.Range("M5").Value = 0.7809334999999999
$ws.Range("N5").Value = 1.561867
$ws.Range("O5").Value = 0.007456697592460336
$ws.Range("P5").Value = 0.004983518592002547
$ws.Range("Q5").Value = 1.916680751679833
$ws.Range("R5").Value = 11.500084510079
$ws.Range("S5").Value = 0.00213072852048838
$ws.Range("T5").Value = 0.001508816807076803
$ws.Range("G6").Value = 2.454345666666666
$ws.Range("H6").Value = 7.363036999999999
$ws.Range("I6").Value = 0.2857469401257222
$ws.Range("J6").Value = 0.3027613480760606
$ws.Range("M6").Value = 3.98709
$ws.Range("N6").Value = 11.96127
$ws.Range("O6").Value = 0.03807049435569441
$ws.Range("P6").Value = 0.03816535686390858
$ws.Range("Q6").Value = 9.78569706411
$ws.Range("R6").Value = 88.07127357699
$ws.Range("S6").Value = 0.01087852727121326
$ws.Range("T6").Value = 0.01155499489392089
$ws.Range("G7").Value = 4.686805000000001
$ws.Range("I7").Value = 0.5456607868665887
$ws.Range("J7").Value = 0.5781514068052169
$ws.Range("M7").Value = 90.43008666666667
$ws.Range("N7").Value = 271.29026
$ws.Range("O7").Value = 0.863466363695901
$ws.Range("P7").Value = 0.8656179140344247
$ws.Range("Q7").Value = 423.8281823397667
$ws.Range("R7").Value = 3814.4536410579
$ws.Range("S7").Value = 0.4711597354471374
$ws.Range("T7").Value = 0.5004582147547999
$ws.Range("G8").Value = 4.686805000000001
$ws.Range("I8").Value = 0.5456607868665887
$ws.Range("J8").Value = 0.5781514068052169
$ws.Range("O8").Value = 0.000602137432244878
$ws.Range("P8").Value = 0.0006036378137891445
$ws.Range("Q8").Value = 0.2955561723733334
$ws.Range("S8").Value = 0.0003285627850805674
$ws.Range("T8").Value = 0.0003489940512430195
$ws.Range("G9").Value = 4.686805000000001
$ws.Range("I9").Value = 0.5456607868665887
$ws.Range("J9").Value = 0.5781514068052169
$ws.Range("M9").Value = 9.467965
$ws.Range("N9").Value = 28.403895
$ws.Range("O9").Value = 0.0904043069236993
$ws.Range("P9").Value = 0.09062957269587499
$ws.Range("Q9").Value = 44.374505701825
$ws.Range("R9").Value = 399.370551316425
$ws.Range("S9").Value = 0.04933008525211435
$ws.Range("T9").Value = 0.0523976149522758
$ws.Range("G10").Value = 4.686805000000001
$ws.Range("I10").Value = 0.5456607868665887
$ws.Range("J10").Value = 0.5781514068052169
$ws.Range("M10").Value = 0.7809334999999999
$ws.Range("N10").Value = 1.561867
$ws.Range("O10").Value = 0.007456697592460336
$ws.Range("P10").Value = 0.004983518592002547
$ws.Range("Q10").Value = 3.6600830324675
$ws.Range("R10").Value = 21.960498194805
$ws.Range("S10").Value = 0.004068827475728104
$ws.Range("T10").Value = 0.002881228284806226
$ws.Range("G11").Value = 4.686805000000001
$ws.Range("I11").Value = 0.5456607868665887
$ws.Range("J11").Value = 0.5781514068052169
$ws.Range("M11").Value = 3.98709
$ws.Range("N11").Value = 11.96127
$ws.Range("O11").Value = 0.03807049435569441
$ws.Range("P11").Value = 0.03816535686390858
$ws.Range("Q11").Value = 18.68671334745
$ws.Range("R11").Value = 168.18042012705
$ws.Range("S11").Value = 0.02077357590652824
$ws.Range("T11").Value = 0.02206535476209189
$ws.Range("G12").Value = 1.4480775
$ws.Range("H12").Value = 2.896155
$ws.Range("I12").Value = 0.1685922730076891
$ws.Range("J12").Value = 0.1190872451187225
$ws.Range("M12").Value = 90.43008666666667
$ws.Range("N12").Value = 271.29026
$ws.Range("O12").Value = 0.863466363695901
$ws.Range("P12").Value = 0.8656179140344247
$ws.Range("Q12").Value = 130.94977382505
$ws.Range("R12").Value = 785.6986429502999
$ws.Range("S12").Value = 0.1455737569211759
$ws.Range("T12").Value = 0.1030840527077748
$ws.Range("G13").Value = 1.4480775
$ws.Range("H13").Value = 2.896155
$ws.Range("I13").Value = 0.1685922730076891
$ws.Range("J13").Value = 0.1190872451187225
$ws.Range("O13").Value = 0.000602137432244878
$ws.Range("P13").Value = 0.0006036378137891445
$ws.Range("Q13").Value = 0.09131769792
$ws.Range("R13").Value = 0.54790618752
$ws.Range("S13").Value = 0.0001015157183651774
$ws.Range("T13").Value = 0.00007188556429363762
$ws.Range("G14").Value = 1.4480775
$ws.Range("H14").Value = 2.896155
$ws.Range("I14").Value = 0.1685922730076891
$ws.Range("J14").Value = 0.1190872451187225
$ws.Range("M14").Value = 9.467965
$ws.Range("N14").Value = 28.403895
$ws.Range("O14").Value = 0.0904043069236993
$ws.Range("P14").Value = 0.09062957269587499
$ws.Range("Q14").Value = 13.7103470872875
$ws.Range("R14").Value = 82.26208252372498
$ws.Range("S14").Value = 0.01524146759395123
$ws.Range("T14").Value = 0.01079282613863874
$ws.Range("G15").Value = 1.4480775
$ws.Range("H15").Value = 2.896155
$ws.Range("I15").Value = 0.1685922730076891
$ws.Range("J15").Value = 0.1190872451187225
$ws.Range("M15").Value = 0.7809334999999999
$ws.Range("N15").Value = 1.561867
$ws.Range("O15").Value = 0.007456697592460336
$ws.Range("P15").Value = 0.004983518592002547
$ws.Range("Q15").Value = 1.13085223034625
$ws.Range("R15").Value = 4.523408921384999
$ws.Range("S15").Value = 0.001257141596243851
$ws.Range("T15").Value = 0.0005934735001195181
$ws.Range("G16").Value = 1.4480775
$ws.Range("H16").Value = 2.896155
$ws.Range("I16").Value = 0.1685922730076891
$ws.Range("J16").Value = 0.1190872451187225
$ws.Range("M16").Value = 3.98709
$ws.Range("N16").Value = 11.96127
$ws.Range("O16").Value = 0.03807049435569441
$ws.Range("P16").Value = 0.03816535686390858
$ws.Range("Q16").Value = 5.773615319475
$ws.Range("R16").Value = 34.64169191685
$ws.Range("S16").Value = 0.006418391177952921
$ws.Range("T16").Value = 0.0045450072078958
